$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.066563650209148
$ws.Range("D2").Value = 1.066085238253936
$ws.Range("E2").Value = 1.070818463736863
$ws.Range("F2").Value = 1.079822210036905
$ws.Range("I2").Value = 1.04687957465538
$ws.Range("J2").Value = 1.071513011321767
$ws.Range("K2").Value = 1.068796736202118
$ws.Range("L2").Value = 1.073517287388646
$ws.Range("M2").Value = 1.082497255433946
$ws.Range("N2").Value = 1.073034682702104
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.068139078074619
$ws.Range("D3").Value = 1.067295218899686
$ws.Range("E3").Value = 1.072202789482804
$ws.Range("F3").Value = 1.081242621272449
$ws.Range("I3").Value = 1.047276844753416
$ws.Range("J3").Value = 1.072741523228966
$ws.Range("K3").Value = 1.069821256481318
$ws.Range("L3").Value = 1.074716644213259
$ws.Range("M3").Value = 1.083734346252194
$ws.Range("N3").Value = 1.074264939237124
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.069156995402744
$ws.Range("D4").Value = 1.068076600552439
$ws.Range("E4").Value = 1.07309736029562
$ws.Range("F4").Value = 1.082160564131754
$ws.Range("I4").Value = 1.047531545720008
$ws.Range("J4").Value = 1.073534558982527
$ws.Range("K4").Value = 1.070482064529698
$ws.Range("L4").Value = 1.075490993015607
$ws.Range("M4").Value = 1.084533137683098
$ws.Range("N4").Value = 1.075059101192417
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.069584579455635
$ws.Range("D5").Value = 1.068404726126916
$ws.Range("E5").Value = 1.07347316165021
$ws.Range("F5").Value = 1.082546196440489
$ws.Range("I5").Value = 1.04763805949215
$ws.Range("J5").Value = 1.073867504087265
$ws.Range("K5").Value = 1.070759364566207
$ws.Range("L5").Value = 1.075816124526414
$ws.Range("M5").Value = 1.084868551137818
$ws.Range("N5").Value = 1.075392519117393
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.069656352509453
$ws.Range("D6").Value = 1.068459798474254
$ws.Range("E6").Value = 1.073536244339181
$ws.Range("F6").Value = 1.082610930103426
$ws.Range("I6").Value = 1.047655910704086
$ws.Range("J6").Value = 1.073923380995274
$ws.Range("K6").Value = 1.070805895074229
$ws.Range("L6").Value = 1.075870691975519
$ws.Range("M6").Value = 1.084924845337418
$ws.Range("N6").Value = 1.07544847537702
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.069162710160498
$ws.Range("D7").Value = 1.068080986418885
$ws.Range("E7").Value = 1.073102382849143
$ws.Range("F7").Value = 1.082165718025963
$ws.Range("I7").Value = 1.047532971169701
$ws.Range("J7").Value = 1.073539009560563
$ws.Range("K7").Value = 1.070485771799173
$ws.Range("L7").Value = 1.075495339021771
$ws.Range("M7").Value = 1.08453762105251
$ws.Range("N7").Value = 1.075063558090784
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.067096386912575
$ws.Range("D8").Value = 1.066494482479499
$ws.Range("E8").Value = 1.071286551401341
$ws.Range("F8").Value = 1.080302488341335
$ws.Range("I8").Value = 1.047014323866277
$ws.Range("J8").Value = 1.071928588501129
$ws.Range("K8").Value = 1.06914342070798
$ws.Range("L8").Value = 1.073922974367812
$ws.Range("M8").Value = 1.082915689358858
$ws.Range("N8").Value = 1.073450850048728
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.063443478731376
$ws.Range("D9").Value = 1.063686693299737
$ws.Range("E9").Value = 1.068077495135307
$ws.Range("F9").Value = 1.077010081093713
$ws.Range("I9").Value = 1.046082226498904
$ws.Range("J9").Value = 1.069076038938988
$ws.Range("K9").Value = 1.066761522139694
$ws.Range("L9").Value = 1.071138864171428
$ws.Range("M9").Value = 1.080044436885525
$ws.Range("N9").Value = 1.070594249539017
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.060999762150454
$ws.Range("D10").Value = 1.061806310102497
$ws.Range("E10").Value = 1.065931446875804
$ws.Range("F10").Value = 1.074808579832095
$ws.Range("I10").Value = 1.045448460490522
$ws.Range("J10").Value = 1.067164007983881
$ws.Range("K10").Value = 1.065162170622696
$ws.Range("L10").Value = 1.069273405398054
$ws.Range("M10").Value = 1.078121018869991
$ws.Range("N10").Value = 1.068679503280641
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.059939486098726
$ws.Range("D11").Value = 1.060989984414047
$ws.Range("E11").Value = 1.065000509691653
$ws.Range("F11").Value = 1.073853657822489
$ws.Range("I11").Value = 1.045171065816826
$ws.Range("J11").Value = 1.066333538948277
$ws.Range("K11").Value = 1.064466852789643
$ws.Range("L11").Value = 1.068463331196344
$ws.Range("M11").Value = 1.077285879842455
$ws.Range("N11").Value = 1.06784785488375
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.059545322109272
$ws.Range("D12").Value = 1.060686441343392
$ws.Range("E12").Value = 1.064654457243028
$ws.Range("F12").Value = 1.073498700330574
$ws.Range("I12").Value = 1.045067579936963
$ws.Range("J12").Value = 1.06602467543129
$ws.Range("K12").Value = 1.064208155568148
$ws.Range("L12").Value = 1.068162077965917
$ws.Range("M12").Value = 1.076975321207296
$ws.Range("N12").Value = 1.067538552745142
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.059629886731578
$ws.Range("D13").Value = 1.060751567104429
$ws.Range("E13").Value = 1.06472869860317
$ws.Range("F13").Value = 1.073574851672586
$ws.Range("I13").Value = 1.045089798382272
$ws.Range("J13").Value = 1.066090945491867
$ws.Range("K13").Value = 1.064263666344824
$ws.Range("L13").Value = 1.068226714016331
$ws.Range("M13").Value = 1.077041953082086
$ws.Range("N13").Value = 1.067604916916806
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.059906911176495
$ws.Range("D14").Value = 1.060964900084943
$ws.Range("E14").Value = 1.064971910257129
$ws.Range("F14").Value = 1.073824322213067
$ws.Range("I14").Value = 1.045162520821441
$ws.Range("J14").Value = 1.066308016201096
$ws.Range("K14").Value = 1.064445477507846
$ws.Range("L14").Value = 1.068438436807059
$ws.Range("M14").Value = 1.07726021616789
$ws.Range("N14").Value = 1.06782229589134
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.060077551012977
$ws.Range("D15").Value = 1.061096298523217
$ws.Range("E15").Value = 1.065121726170926
$ws.Range("F15").Value = 1.073977994961311
$ws.Range("I15").Value = 1.045207267910769
$ws.Range("J15").Value = 1.066441708674532
$ws.Range("K15").Value = 1.064557440829529
$ws.Range("L15").Value = 1.06856883886546
$ws.Range("M15").Value = 1.077394648552125
$ws.Range("N15").Value = 1.067956178223423
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.061070083425307
$ws.Range("D16").Value = 1.061860441975018
$ws.Range("E16").Value = 1.065993193954658
$ws.Range("F16").Value = 1.074871919265396
$ws.Range("I16").Value = 1.045466807438392
$ws.Range("J16").Value = 1.067219069161906
$ws.Range("K16").Value = 1.065208257295684
$ws.Range("L16").Value = 1.069327117889296
$ws.Range("M16").Value = 1.078176395502955
$ws.Range("N16").Value = 1.068734642651855
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.061692095210015
$ws.Range("D17").Value = 1.062339200023303
$ws.Range("E17").Value = 1.066539386220305
$ws.Range("F17").Value = 1.07543220504439
$ws.Range("I17").Value = 1.045628812539164
$ws.Range("J17").Value = 1.06770599954731
$ws.Range("K17").Value = 1.065615746363558
$ws.Range("L17").Value = 1.069802140483109
$ws.Range("M17").Value = 1.078666147446314
$ws.Range("N17").Value = 1.069222264534257
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.062054699274124
$ws.Range("D18").Value = 1.062618248623853
$ws.Range("E18").Value = 1.066857808972318
$ws.Range("F18").Value = 1.07575885088111
$ws.Range("I18").Value = 1.045723021032231
$ws.Range("J18").Value = 1.067989772575552
$ws.Range("K18").Value = 1.065853159266049
$ws.Range("L18").Value = 1.070078989946151
$ws.Range("M18").Value = 1.078951591369044
$ws.Range("N18").Value = 1.069506440552749
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.062178303442791
$ws.Range("D19").Value = 1.062713362778476
$ws.Range("E19").Value = 1.066966355669457
$ws.Range("F19").Value = 1.075870201875295
$ws.Range("I19").Value = 1.045755095214499
$ws.Range("J19").Value = 1.068086490481154
$ws.Range("K19").Value = 1.06593406556604
$ws.Range("L19").Value = 1.070173350823663
$ws.Range("M19").Value = 1.079048883257651
$ws.Range("N19").Value = 1.06960329580887
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.061625380498045
$ws.Range("D20").Value = 1.062287854848319
$ws.Range("E20").Value = 1.066480801768462
$ws.Range("F20").Value = 1.075372108225569
$ws.Range("I20").Value = 1.045611460564179
$ws.Range("J20").Value = 1.067653781945885
$ws.Range("K20").Value = 1.065572054469121
$ws.Range("L20").Value = 1.069751198191834
$ws.Range("M20").Value = 1.078613624516771
$ws.Range("N20").Value = 1.069169972777849
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.059825343580683
$ws.Range("D21").Value = 1.060902087820328
$ws.Range("E21").Value = 1.064900297753655
$ws.Range("F21").Value = 1.07375086651358
$ws.Range("I21").Value = 1.045141118298291
$ws.Range("J21").Value = 1.066244105124878
$ws.Range("K21").Value = 1.064391950441208
$ws.Range("L21").Value = 1.068376099597664
$ws.Range("M21").Value = 1.077195952872507
$ws.Range("N21").Value = 1.067758294054063
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.058691673088401
$ws.Range("D22").Value = 1.060028926819788
$ws.Range("E22").Value = 1.06390505853464
$ws.Range("F22").Value = 1.072730036953697
$ws.Range("I22").Value = 1.044842795331121
$ws.Range("J22").Value = 1.065355523117308
$ws.Range("K22").Value = 1.063647508831063
$ws.Range("L22").Value = 1.067509459456978
$ws.Range("M22").Value = 1.076302573060482
$ws.Range("N22").Value = 1.066868450158121
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.059292838274059
$ws.Range("D23").Value = 1.060491985916669
$ws.Range("E23").Value = 1.064432799804379
$ws.Range("F23").Value = 1.073271341968407
$ws.Range("I23").Value = 1.045001189439153
$ws.Range("J23").Value = 1.065826794266327
$ws.Range("K23").Value = 1.064042387044356
$ws.Range("L23").Value = 1.067969079615733
$ws.Range("M23").Value = 1.076776365802369
$ws.Range("N23").Value = 1.067340390566227
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.061655526637965
$ws.Range("D24").Value = 1.062311056149868
$ws.Range("E24").Value = 1.066507274061188
$ws.Range("F24").Value = 1.075399263883889
$ws.Range("I24").Value = 1.04561930205965
$ws.Range("J24").Value = 1.067677377590228
$ws.Range("K24").Value = 1.065591797785604
$ws.Range("L24").Value = 1.069774217510292
$ws.Range("M24").Value = 1.078637358048129
$ws.Range("N24").Value = 1.069193601930713
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.064389292914627
$ws.Range("D25").Value = 1.064414051100424
$ws.Range("E25").Value = 1.068908260166223
$ws.Range("F25").Value = 1.0778623741588
$ws.Range("I25").Value = 1.046325364956648
$ws.Range("J25").Value = 1.069815283349584
$ws.Range("K25").Value = 1.067379288905915
$ws.Range("L25").Value = 1.071860250021428
$ws.Range("M25").Value = 1.080788327219139
$ws.Range("N25").Value = 1.071334543761468

Write-Output "Applied 264 cell updates"